$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents (keep formatting) of the rows that no longer
# contain the "Ressources Humaines"/K'IAM, Yoobic and Risorse Umane/SMARTRH
# entries. Excel's ClearContents removes the cell value but leaves the
# cell style untouched, matching the target XML where only the t="s"
# attribute and <v> element disappear while s="..." is kept.
$rowsToClear = @(4, 18, 19, 20, 32, 33, 39, 40, 44)

foreach ($r in $rowsToClear) {
    $ws.Range("A$r`:D$r").ClearContents()
}

# Reflect the final selection left by the author (cell D44).
$ws.Range("D44").Select()
